$d = $word.ActiveDocument

# 1. "Week of code 23" -> "Week of code 28"
$d.Content.Find.Execute("Week of code 23", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Week of code 28", 2)

# 2. "720/10489" -> "385/10432" (collapses three runs into the text of the first run,
#    the other two runs' text is cleared)
$d.Content.Find.Execute("720/10489", $true, $false, $false, $false, $false,
                         $true, 1, $false, "385/10432", 2)

# 3. Merge "Week of code 27 Rank: 246/" + "7941" into a single run's text
$d.Content.Find.Execute("Week of code 27 Rank: 246/7941", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Week of code 27 Rank: 246/7941", 2)
